$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 348.129669
$ws.Range("H2").Value = 1044.389007
$ws.Range("I2").Value = 0.6448154080547559
$ws.Range("J2").Value = 0.644815408054756
$ws.Range("M2").Value = 281.0920463333333
$ws.Range("N2").Value = 843.2761389999999
$ws.Range("O2").Value = 0.8291026083535286
$ws.Range("P2").Value = 0.8291026083535286
$ws.Range("Q2").Value = 97856.48104855599
$ws.Range("R2").Value = 880708.3294370039
$ws.Range("S2").Value = 0.534618136724743
$ws.Range("T2").Value = 0.5346181367247431

# Row 3
$ws.Range("G3").Value = 348.129669
$ws.Range("H3").Value = 1044.389007
$ws.Range("I3").Value = 0.6448154080547559
$ws.Range("J3").Value = 0.644815408054756
$ws.Range("O3").Value = 0.001324719879221983
$ws.Range("P3").Value = 0.001324719879221983
$ws.Range("Q3").Value = 156.352693200618
$ws.Range("R3").Value = 1407.174238805562
$ws.Range("S3").Value = 0.0008541997894787702
$ws.Range("T3").Value = 0.0008541997894787704

# Row 4
$ws.Range("G4").Value = 348.129669
$ws.Range("H4").Value = 1044.389007
$ws.Range("I4").Value = 0.6448154080547559
$ws.Range("J4").Value = 0.644815408054756
$ws.Range("M4").Value = 4.452417
$ws.Range("N4").Value = 13.357251
$ws.Range("O4").Value = 0.01313274635953239
$ws.Range("P4").Value = 0.01313274635953239
$ws.Range("Q4").Value = 1550.018456459973
$ws.Range("R4").Value = 13950.16610813976
$ws.Range("S4").Value = 0.008468197202701486
$ws.Range("T4").Value = 0.008468197202701489

# Row 5
$ws.Range("G5").Value = 348.129669
$ws.Range("H5").Value = 1044.389007
$ws.Range("I5").Value = 0.6448154080547559
$ws.Range("J5").Value = 0.644815408054756
$ws.Range("M5").Value = 53.03808999999999
$ws.Range("N5").Value = 159.11427
$ws.Range("O5").Value = 0.156439925407717
$ws.Range("P5").Value = 0.156439925407717
$ws.Range("Q5").Value = 18464.1327160922
$ws.Range("R5").Value = 166177.1944448299
$ws.Range("S5").Value = 0.1008748743378326
$ws.Range("T5").Value = 0.1008748743378326

# Row 6
$ws.Range("I6").Value = 0.1713626578780604
$ws.Range("J6").Value = 0.1713626578780604
$ws.Range("M6").Value = 281.0920463333333
$ws.Range("N6").Value = 843.2761389999999
$ws.Range("O6").Value = 0.8291026083535286
$ws.Range("P6").Value = 0.8291026083535286
$ws.Range("Q6").Value = 26005.80952874908
$ws.Range("R6").Value = 234052.2857587417
$ws.Range("S6").Value = 0.1420772266210932
$ws.Range("T6").Value = 0.1420772266210932

# Row 7
$ws.Range("I7").Value = 0.1713626578780604
$ws.Range("J7").Value = 0.1713626578780604
$ws.Range("O7").Value = 0.001324719879221983
$ws.Range("P7").Value = 0.001324719879221983
$ws.Range("S7").Value = 0.0002270075194473822
$ws.Range("T7").Value = 0.0002270075194473822

# Row 8
$ws.Range("I8").Value = 0.1713626578780604
$ws.Range("J8").Value = 0.1713626578780604
$ws.Range("M8").Value = 4.452417
$ws.Range("N8").Value = 13.357251
$ws.Range("O8").Value = 0.01313274635953239
$ws.Range("P8").Value = 0.01313274635953239
$ws.Range("Q8").Value = 411.924527765742
$ws.Range("R8").Value = 3707.320749891678
$ws.Range("S8").Value = 0.002250462321407892
$ws.Range("T8").Value = 0.002250462321407892

# Row 9
$ws.Range("I9").Value = 0.1713626578780604
$ws.Range("J9").Value = 0.1713626578780604
$ws.Range("M9").Value = 53.03808999999999
$ws.Range("N9").Value = 159.11427
$ws.Range("O9").Value = 0.156439925407717
$ws.Range("P9").Value = 0.156439925407717
$ws.Range("Q9").Value = 4906.928119456673
$ws.Range("R9").Value = 44162.35307511006
$ws.Range("S9").Value = 0.02680796141611189
$ws.Range("T9").Value = 0.02680796141611189

# Row 10
$ws.Range("G10").Value = 98.04896266666667
$ws.Range("H10").Value = 294.146888
$ws.Range("I10").Value = 0.1816090023377243
$ws.Range("J10").Value = 0.1816090023377243
$ws.Range("M10").Value = 281.0920463333333
$ws.Range("N10").Value = 843.2761389999999
$ws.Range("O10").Value = 0.8291026083535286
$ws.Range("P10").Value = 0.8291026083535286
$ws.Range("Q10").Value = 27560.78355683394
$ws.Range("R10").Value = 248047.0520115054
$ws.Range("S10").Value = 0.1505724975386893
$ws.Range("T10").Value = 0.1505724975386893

# Row 11
$ws.Range("G11").Value = 98.04896266666667
$ws.Range("H11").Value = 294.146888
$ws.Range("I11").Value = 0.1816090023377243
$ws.Range("J11").Value = 0.1816090023377243
$ws.Range("O11").Value = 0.001324719879221983
$ws.Range("P11").Value = 0.001324719879221983
$ws.Range("Q11").Value = 44.03594621077867
$ws.Range("R11").Value = 396.323515897008
$ws.Range("S11").Value = 0.000240581055642455
$ws.Range("T11").Value = 0.0002405810556424551

# Row 12
$ws.Range("G12").Value = 98.04896266666667
$ws.Range("H12").Value = 294.146888
$ws.Range("I12").Value = 0.1816090023377243
$ws.Range("J12").Value = 0.1816090023377243
$ws.Range("M12").Value = 4.452417
$ws.Range("N12").Value = 13.357251
$ws.Range("O12").Value = 0.01313274635953239
$ws.Range("P12").Value = 0.01313274635953239
$ws.Range("Q12").Value = 436.554868209432
$ws.Range("R12").Value = 3928.993813884888
$ws.Range("S12").Value = 0.002385024964309058
$ws.Range("T12").Value = 0.002385024964309058

# Row 13
$ws.Range("G13").Value = 98.04896266666667
$ws.Range("H13").Value = 294.146888
$ws.Range("I13").Value = 0.1816090023377243
$ws.Range("J13").Value = 0.1816090023377243
$ws.Range("M13").Value = 53.03808999999999
$ws.Range("N13").Value = 159.11427
$ws.Range("O13").Value = 0.156439925407717
$ws.Range("P13").Value = 0.156439925407717
$ws.Range("Q13").Value = 5200.329706321306
$ws.Range("R13").Value = 46802.96735689175
$ws.Range("S13").Value = 0.02841089877908349
$ws.Range("T13").Value = 0.0284108987790835

# Row 14
$ws.Range("G14").Value = 1.194740666666667
$ws.Range("H14").Value = 3.584222
$ws.Range("I14").Value = 0.002212931729459341
$ws.Range("J14").Value = 0.002212931729459341
$ws.Range("M14").Value = 281.0920463333333
$ws.Range("N14").Value = 843.2761389999999
$ws.Range("O14").Value = 0.8291026083535286
$ws.Range("P14").Value = 0.8291026083535286
$ws.Range("Q14").Value = 335.8320988309842
$ws.Range("R14").Value = 3022.488889478858
$ws.Range("S14").Value = 0.001834747469003024
$ws.Range("T14").Value = 0.001834747469003025

# Row 15
$ws.Range("G15").Value = 1.194740666666667
$ws.Range("H15").Value = 3.584222
$ws.Range("I15").Value = 0.002212931729459341
$ws.Range("J15").Value = 0.002212931729459341
$ws.Range("O15").Value = 0.001324719879221983
$ws.Range("P15").Value = 0.001324719879221983
$ws.Range("Q15").Value = 0.5365843176946667
$ws.Range("R15").Value = 4.829258859252
$ws.Range("S15").Value = 0.000002931514653375872
$ws.Range("T15").Value = 0.000002931514653375873

# Row 16
$ws.Range("G16").Value = 1.194740666666667
$ws.Range("H16").Value = 3.584222
$ws.Range("I16").Value = 0.002212931729459341
$ws.Range("J16").Value = 0.002212931729459341
$ws.Range("M16").Value = 4.452417
$ws.Range("N16").Value = 13.357251
$ws.Range("O16").Value = 0.01313274635953239
$ws.Range("P16").Value = 0.01313274635953239
$ws.Range("Q16").Value = 5.319483654858
$ws.Range("R16").Value = 47.875352893722
$ws.Range("S16").Value = 0.00002906187111395087
$ws.Range("T16").Value = 0.00002906187111395088

# Row 17
$ws.Range("G17").Value = 1.194740666666667
$ws.Range("H17").Value = 3.584222
$ws.Range("I17").Value = 0.002212931729459341
$ws.Range("J17").Value = 0.002212931729459341
$ws.Range("M17").Value = 53.03808999999999
$ws.Range("N17").Value = 159.11427
$ws.Range("O17").Value = 0.156439925407717
$ws.Range("P17").Value = 0.156439925407717
$ws.Range("Q17").Value = 63.36676300532665
$ws.Range("R17").Value = 570.30086704794
$ws.Range("S17").Value = 0.0003461908746889894
$ws.Range("T17").Value = 0.0003461908746889895
